$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The report gained two additional threat rows. Concretely:
#   - a new row was inserted right before the old "31-JAN-26" row (old row 3)
#   - the (now shifted) "11-FEB-26" row had its Our Fare / Fare Dif updated
#   - a new row was inserted right before the old "28-MAR-26" row (old row 6)
# This pushes the former rows 3-7 down to rows 4,5,6,8,9 and adds two
# brand-new rows (3 and 7).
# ---------------------------------------------------------------------------

# --- Insert new row 3 and fill it in -----------------------------------
$ws.Rows.Item(3).Insert()

# Force column A to be text so the "DD-MMM-YY" string isn't auto-converted
# into a date serial number, then set all the values for the new row.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "22-JAN-26"
$ws.Range("B3").Value = "SM-443"
$ws.Range("C3").Value = "Nile Air NP-143"
$ws.Range("D3").Value = 13655
$ws.Range("E3").Value = 15695
$ws.Range("F3").Value = -2040
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "LOW THREAT"
$ws.Range("K3").Value = "EGP"

# Re-apply the standard row formatting (border/fill/font/alignment) from the
# row right below, which restores the shared cell style used throughout the
# table (this only touches formats, not the values set above).
$ws.Range("A4:K4").Copy()
$ws.Range("A3:K3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 5 (shifted old row 4): Our Fare / Fare Dif updated ----------------
$ws.Range("E5").Value = 7923
$ws.Range("F5").Value = -725

# --- Insert new row 7 and fill it in -----------------------------------
$ws.Rows.Item(7).Insert()

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "25-FEB-26"
$ws.Range("B7").Value = "SM-447"
$ws.Range("C7").Value = "Air Arabia Egypt E5-513"
$ws.Range("D7").Value = 7198
$ws.Range("E7").Value = 7495
$ws.Range("F7").Value = -297
$ws.Range("G7").Value = 30
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = "LOW THREAT"
$ws.Range("K7").Value = "EGP"

$ws.Range("A6:K6").Copy()
$ws.Range("A7:K7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
